# The deck's slide master currently uses the "Integral" theme palette
# (ppt/theme/theme1.xml). The commit swaps the presentation's effective
# theme colours to the stock "Office Theme" palette that was previously
# only sitting, unused, in ppt/theme/theme2.xml.
#
# PowerPoint's ColorScheme object maps the 12 theme colour slots as:
#   1 dk1  2 lt1  3 dk2  4 lt2  5 accent1 6 accent2
#   7 accent3 8 accent4 9 accent5 10 accent6 11 hlink 12 folHlink
# .RGB uses the VBA long-colour encoding (R + G*256 + B*65536), so each
# target hex colour below is pre-converted to that decimal form.

$p  = $ppt.ActivePresentation
$cs = $p.SlideMaster.ColorScheme

$cs.Colors(1).RGB  = 0          # dk1      000000 (unchanged)
$cs.Colors(2).RGB  = 16777215   # lt1      FFFFFF (unchanged)
$cs.Colors(3).RGB  = 6968388    # dk2      44546A
$cs.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$cs.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$cs.Colors(6).RGB  = 3243501    # accent2  ED7D31
$cs.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$cs.Colors(8).RGB  = 49407      # accent4  FFC000
$cs.Colors(9).RGB  = 12874308   # accent5  4472C4
$cs.Colors(10).RGB = 4697456    # accent6  70AD47
$cs.Colors(11).RGB = 12673797   # hlink    0563C1
$cs.Colors(12).RGB = 7491477    # folHlink 954F72
